$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update map tile values (texture/lighting tweaks on the Pac-Man map)
$ws.Range("I5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("C19").Value = 2
$ws.Range("E19").Value = 2
$ws.Range("G19").Value = 2
$ws.Range("K19").Value = 2

# Select the whole sheet (matches selection sqref="A1:XFD1048576")
$ws.Cells.Select()
